# Weekly fruit/vegetable price update: insert two new daily records at the
# top of the data block (new rows 450-451), pushing the existing rows
# 450-481 down to 452-483.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 450.
$ws.Rows("450:451").Insert()

# --- New row 450 ---------------------------------------------------------
$ws.Cells.Item(450, 1).Value = 9
$ws.Cells.Item(450, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(450, 3).Value = "Metropolitana"
$ws.Cells.Item(450, 4).Value = 44585
$ws.Cells.Item(450, 5).Value = 13
$ws.Cells.Item(450, 6).Value = 100112040
$ws.Cells.Item(450, 7).Value = "Cilantro"
$ws.Cells.Item(450, 8).Value = "Sin especificar"
$ws.Cells.Item(450, 9).Value = "Primera"
$ws.Cells.Item(450, 10).Value = 43
$ws.Cells.Item(450, 11).Value = 8000
$ws.Cells.Item(450, 12).Value = 8000
$ws.Cells.Item(450, 13).Value = 8000
$ws.Cells.Item(450, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(450, 15).Value = "Región Metropolitana"
$ws.Cells.Item(450, 16).Value = 222
$ws.Cells.Item(450, 17).Value = 36
$ws.Cells.Item(450, 18).Value = "Hortaliza"

# --- New row 451 ---------------------------------------------------------
$ws.Cells.Item(451, 1).Value = 9
$ws.Cells.Item(451, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(451, 3).Value = "Metropolitana"
$ws.Cells.Item(451, 4).Value = 44585
$ws.Cells.Item(451, 5).Value = 13
$ws.Cells.Item(451, 6).Value = 100112040
$ws.Cells.Item(451, 7).Value = "Cilantro"
$ws.Cells.Item(451, 8).Value = "Sin especificar"
$ws.Cells.Item(451, 9).Value = "Primera"
$ws.Cells.Item(451, 10).Value = 79
$ws.Cells.Item(451, 11).Value = 16000
$ws.Cells.Item(451, 12).Value = 18000
$ws.Cells.Item(451, 13).Value = 16987
$ws.Cells.Item(451, 14).Value = "$/docena de atados"
$ws.Cells.Item(451, 15).Value = "Región Metropolitana"
$ws.Cells.Item(451, 16).Value = 5662
$ws.Cells.Item(451, 17).Value = 3
$ws.Cells.Item(451, 18).Value = "Hortaliza"
